$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.810.66"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "1.888.84"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7674"
$ws.Range("E5").Value = "  -4.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.23"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3120"
$ws.Range("E8").Value = "  -4.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.25"
$ws.Range("E9").Value = "  -7.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07205"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08073"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7644"
$ws.Range("E12").Value = "  -4.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.521"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "1.855.95"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.20"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.137"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "29.819.45"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.02"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007760"
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.151.74"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.133"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1546"
$ws.Range("E25").Value = "  -5.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.389"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.18"
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.036"
$ws.Range("E29").Value = "  -5.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.440"
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.548"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.080"
$ws.Range("E33").Value = "  -1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05498"
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.255"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7469"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9986"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.631"
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01919"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.779"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("D41").Value = "1.148.30"
$ws.Range("E41").Value = "  +11.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4407"
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "73.30"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.876"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8501"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.50"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.880"
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.879"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.441"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.983"
$ws.Range("E51").Value = "  +8.98%  "
